# Append a new job posting row (JD_019 / Senior QA Analyst) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the first empty row below the existing data (row 20, since data currently runs rows 1-19).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Job description text - reuses the same "Junior RPA Developer" description already
# present in the sheet (e.g. the Job_Description for JD_003 / JD_010 / JD_018).
$desc = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions." + [char]10 + "Collaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"

$ws.Cells.Item($newRow, 1).Value = "JD_019"
$ws.Cells.Item($newRow, 2).Value = "Senior QA Analyst"
$ws.Cells.Item($newRow, 3).Value = $desc
$ws.Cells.Item($newRow, 4).Value = 1
$ws.Cells.Item($newRow, 5).Value = 4

# Re-fit the row height so saving doesn't leave a stray custom row height
# caused by the embedded newline in the description text.
$ws.Rows.Item($newRow).AutoFit()

$wb.Save()
